# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, mirroring the style already used by the other header cells (e.g. AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Team record is the same for every player row (2-51): 97 wins, 65 losses, 0 ties.
$ws.Range("AD2:AD51").Value = 97
$ws.Range("AE2:AE51").Value = 65
$ws.Range("AF2:AF51").Value = 0
